$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: write each date as a formula producing text, then convert to values ---
# (A plain string like "02-08-2021" would be auto-parsed by Excel as a date serial;
#  wrapping it in a formula that evaluates to text and then collapsing the formula
#  via Copy + PasteSpecial(xlPasteValues) keeps it as literal text, matching the source file.)
$ws.Range("A150").Formula = "=""02-08-2021"""
$ws.Range("A151").Formula = "=""03-08-2021"""
$ws.Range("A152").Formula = "=""04-08-2021"""
$ws.Range("A153").Formula = "=""05-08-2021"""
$ws.Range("A154").Formula = "=""06-08-2021"""
$ws.Range("A155").Formula = "=""09-08-2021"""
$ws.Range("A156").Formula = "=""10-08-2021"""
$ws.Range("A157").Formula = "=""11-08-2021"""
$ws.Range("A158").Formula = "=""12-08-2021"""
$ws.Range("A159").Formula = "=""13-08-2021"""
$ws.Range("A160").Formula = "=""16-08-2021"""
$ws.Range("A161").Formula = "=""17-08-2021"""
$ws.Range("A162").Formula = "=""18-08-2021"""
$ws.Range("A163").Formula = "=""19-08-2021"""
$ws.Range("A164").Formula = "=""20-08-2021"""
$ws.Range("A165").Formula = "=""23-08-2021"""
$ws.Range("A166").Formula = "=""24-08-2021"""
$ws.Range("A167").Formula = "=""25-08-2021"""
$ws.Range("A168").Formula = "=""26-08-2021"""
$ws.Range("A169").Formula = "=""27-08-2021"""
$ws.Range("A170").Formula = "=""30-08-2021"""
$ws.Range("A171").Formula = "=""31-08-2021"""
$ws.Range("A172").Formula = "=""01-09-2021"""
$ws.Range("A173").Formula = "=""02-09-2021"""
$ws.Range("A174").Formula = "=""03-09-2021"""

$dateRange = $ws.Range("A150:A174")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# --- Columns B, C, D: numeric values ---
$ws.Range("B150").Value = 3.43
$ws.Range("C150").Value = 3.15
$ws.Range("D150").Value = 3.18
$ws.Range("B151").Value = 3.43
$ws.Range("C151").Value = 3.17
$ws.Range("D151").Value = 3.16
$ws.Range("B152").Value = 3.42
$ws.Range("C152").Value = 3.2
$ws.Range("D152").Value = 3.19
$ws.Range("B153").Value = 3.34
$ws.Range("C153").Value = 3.19
$ws.Range("D153").Value = 3.18
$ws.Range("B154").Value = 3.23
$ws.Range("C154").Value = 3.2
$ws.Range("D154").Value = 3.19
$ws.Range("B155").Value = 3.29
$ws.Range("C155").Value = 3.19
$ws.Range("D155").Value = 3.2
$ws.Range("B156").Value = 3.34
$ws.Range("C156").Value = 3.17
$ws.Range("D156").Value = 3.21
$ws.Range("B157").Value = 3.48
$ws.Range("C157").Value = 3.24
$ws.Range("D157").Value = 3.29
$ws.Range("B158").Value = 3.45
$ws.Range("C158").Value = 3.29
$ws.Range("D158").Value = 3.31
$ws.Range("B159").Value = 3.53
$ws.Range("C159").Value = 3.34
$ws.Range("D159").Value = 3.28
$ws.Range("B160").Value = 3.69
$ws.Range("C160").Value = 3.41
$ws.Range("D160").Value = 3.26
$ws.Range("B161").Value = 3.6
$ws.Range("C161").Value = 3.32
$ws.Range("D161").Value = 3.23
$ws.Range("B162").Value = 3.55
$ws.Range("C162").Value = 3.25
$ws.Range("D162").Value = 3.16
$ws.Range("B163").Value = 3.53
$ws.Range("C163").Value = 3.27
$ws.Range("D163").Value = 3.19
$ws.Range("B164").Value = 3.61
$ws.Range("C164").Value = 3.3
$ws.Range("D164").Value = 3.18
$ws.Range("B165").Value = 3.66
$ws.Range("C165").Value = 3.32
$ws.Range("D165").Value = 3.23
$ws.Range("B166").Value = 3.65
$ws.Range("C166").Value = 3.29
$ws.Range("D166").Value = 3.22
$ws.Range("B167").Value = 3.69
$ws.Range("C167").Value = 3.29
$ws.Range("D167").Value = 3.21
$ws.Range("B168").Value = 3.74
$ws.Range("C168").Value = 3.31
$ws.Range("D168").Value = 3.17
$ws.Range("B169").Value = 3.72
$ws.Range("C169").Value = 3.27
$ws.Range("D169").Value = 3.14
$ws.Range("B170").Value = 3.75
$ws.Range("C170").Value = 3.26
$ws.Range("D170").Value = 3.15
$ws.Range("B171").Value = 3.59
$ws.Range("C171").Value = 3.23
$ws.Range("D171").Value = 3.18
$ws.Range("B172").Value = 3.74
$ws.Range("C172").Value = 3.13
$ws.Range("D172").Value = 3.1
$ws.Range("B173").Value = 3.78
$ws.Range("C173").Value = 3.15
$ws.Range("D173").Value = 3.16
$ws.Range("B174").Value = 3.72
$ws.Range("C174").Value = 3.17
$ws.Range("D174").Value = 3.22

